$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.264.54"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").Value = "1.593.10"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'212.82"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.57%  "

$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("D10").Value = "'18.96"
$ws.Range("E10").Value = "  -2.10%  "

$ws.Range("D11").Value = "'0.0850"
$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").Value = "1.817.03"
$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("D13").Value = "1.593.28"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("D15").Value = "'0.508"

$ws.Range("D16").Value = "'63.94"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").Value = "26.261.23"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.40"
$ws.Range("E19").Value = "  -0.95%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'215.31"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").Value = "'4.30"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  -2.62%  "

$ws.Range("D25").Value = "'144.66"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D27").Value = "'6.97"
$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("E30").Value = "  -1.09%  "

$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("D33").Value = "1.426.84"
$ws.Range("E33").Value = "  +6.63%  "

$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  -1.00%  "

$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("D37").Value = "'0.565"
$ws.Range("E37").Value = "  -4.65%  "

$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("D40").Value = "'5.77"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("E42").Value = "  +0.88%  "

$ws.Range("E43").Value = "  -8.64%  "

$ws.Range("D44").Value = "'0.759"
$ws.Range("E44").Value = "  -0.45%  "

$ws.Range("D45").Value = "1.730.13"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("D46").Value = "'60.84"

$ws.Range("D47").Value = "'86.61"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.47"
$ws.Range("E48").Value = "  -1.87%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0500"
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0952"
$ws.Range("E50").Value = "  -2.91%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  +0.03%  "
